$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of profit data appended after a run on 2026-01-04.
# Force column A to stay plain text (not auto-converted to a date serial),
# then strip the format Excel implicitly applied for the "looks like a
# date" text entry so the new row matches the unstyled data rows above it.
$ws.Range("A41").Value = "'01/04/2026"
$ws.Range("A41").ClearFormats()

$ws.Range("B41").Value = 13302.26
$ws.Range("C41").Value = 0.2046055683399551
$ws.Range("D41").Value = 0.7953944316600449
$ws.Range("E41").Value = -97.15000000000001
$ws.Range("F41").Value = -16.61
$ws.Range("G41").Value = -20061.21
$ws.Range("H41").Value = -65.47
$ws.Range("I41").Value = -381.13
$ws.Range("J41").Value = -12.28
